$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new daily price-record row is inserted before row 16. The easiest way to
# reproduce the resulting row (which keeps every column identical to the old
# row 16 except for a handful of edited fields) is to duplicate row 16 in
# place - copy it, insert the copy above itself (pushing the old row 16 and
# everything below it down by one), and then edit the fields that changed
# for the new record.
$ws.Rows(16).Copy()
$ws.Rows(16).Insert()

# Update the new row 16 with the new record's data.
$ws.Range("D16").Value = 44910
$ws.Range("N16").Value = 3000
$ws.Range("O16").Value = 3000
$ws.Range("P16").Value = 3000
$ws.Range("R16").Value = "Provincia de Curicó"
$ws.Range("S16").Value = 1500
